$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 100084927
$ws.Range("B2").Value = 88856
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 2008
$ws.Range("F2").Value = "Fyrflikig jordstjärna"
$ws.Range("G2").Value = "Geastrum quadrifidum"
$ws.Range("H2").Value = "Pers.:Pers."
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "6"
$ws.Range("J2").Value = "fruktkroppar"
$ws.Range("K2").Value = "'"
$ws.Range("P2").Value = "Skjulstagatan 17, Eskilstuna, Srm"
$ws.Range("Q2").Value = 583689.2563185043
$ws.Range("R2").Value = 6579215.892053389
$ws.Range("S2").Value = 25
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2022-04-19"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2022-04-19"
$ws.Range("AI2").ClearContents()
$ws.Range("AW2").Value = "Dan Åman"
$ws.Range("AX2").Value = "Dan Åman"

# Row 3
$ws.Range("A3").Value = 83948609
$ws.Range("B3").Value = 98520
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("K3").Value = "blomning"
$ws.Range("Q3").Value = 583475.6238590981
$ws.Range("R3").Value = 6579331.079241654
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2020-03-22"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2020-03-22"
$ws.Range("AI3").Value = "Blandskog, kalkberg"

# Row 4
$ws.Range("A4").Value = 86813759
$ws.Range("B4").Value = 103813
$ws.Range("D4").Value = "EN"
$ws.Range("E4").Value = 220785
$ws.Range("F4").Value = "Ask"
$ws.Range("G4").Value = "Fraxinus excelsior"
$ws.Range("H4").Value = "L."

# Row 5
$ws.Range("A5").Value = 86813766
$ws.Range("B5").Value = 98520
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 222498
$ws.Range("F5").Value = "Blåsippa"
$ws.Range("G5").Value = "Hepatica nobilis"
$ws.Range("H5").Value = "Schreb."
$ws.Range("I5").Value = "'"
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("P5").Value = "Vilsta industriområde, V om, Srm"
$ws.Range("Q5").Value = 583486.5729995462
$ws.Range("R5").Value = 6579251.726638615
$ws.Range("S5").Value = 10
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2020-07-09"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2020-07-09"
$ws.Range("AI5").Value = "Blandskog"
$ws.Range("AW5").Value = "Håkan Gustafson"
$ws.Range("AX5").Value = "Håkan Gustafson"
